$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 6-9 (data reduced from 8 data rows to 4 data rows)
$ws.Rows("6:9").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Timp3"
$ws.Range("C2").Value = "Agtr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 69.42974966666667
$ws.Range("H2").Value = 208.289249
$ws.Range("I2").Value = 0.6762611189535584
$ws.Range("J2").Value = 0.6762611189535584
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9721403333333333
$ws.Range("N2").Value = 2.916421
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 67.49545998420322
$ws.Range("R2").Value = 607.459139857829
$ws.Range("S2").Value = 0.6762611189535584
$ws.Range("T2").Value = 0.6762611189535584

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Timp3"
$ws.Range("C3").Value = "Agtr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.70856733333333
$ws.Range("H3").Value = 47.125702
$ws.Range("I3").Value = 0.1530049204123442
$ws.Range("J3").Value = 0.1530049204123442
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9721403333333333
$ws.Range("N3").Value = 2.916421
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 15.27093188361578
$ws.Range("R3").Value = 137.438386952542
$ws.Range("S3").Value = 0.1530049204123442
$ws.Range("T3").Value = 0.1530049204123442

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Timp3"
$ws.Range("C4").Value = "Agtr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.31398133333333
$ws.Range("H4").Value = 51.941944
$ws.Range("I4").Value = 0.1686420078746507
$ws.Range("J4").Value = 0.1686420078746507
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9721403333333333
$ws.Range("N4").Value = 2.916421
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 16.83161958471378
$ws.Range("R4").Value = 151.484576262424
$ws.Range("S4").Value = 0.1686420078746507
$ws.Range("T4").Value = 0.1686420078746507

# Row 5
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("B5").Value = "Timp3"
$ws.Range("C5").Value = "Agtr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2147746666666667
$ws.Range("H5").Value = 0.644324
$ws.Range("I5").Value = 0.002091952759446708
$ws.Range("J5").Value = 0.002091952759446708
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9721403333333333
$ws.Range("N5").Value = 2.916421
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.2087911160448889
$ws.Range("R5").Value = 1.879120044404
$ws.Range("S5").Value = 0.002091952759446708
$ws.Range("T5").Value = 0.002091952759446708
